$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 44201.9993055556
$ws.Range("A3").NumberFormat = "mm/dd/yy\ hh:mm\ AM/PM"

$ws.Range("B3").Value = "Weekly Adventures: Verilog 1"
$ws.Range("C3").Value = "Weekly Adventures"
$ws.Range("D3").Value = "Weekly_Adventures/2021_05_01.md"

$ws.Range("A3").Select()
